# Update RegularMenu assets: reorder/rename menu items and their prices.
#
# Strategy: PowerPoint's TextRange.Text setter here does a naive
# common-prefix diff against the existing text, which can split a
# paragraph's text into extra runs when the new value shares a prefix
# with the old one (e.g. "$2.50" -> "$2.00"). To guarantee a single,
# clean run per paragraph (matching the target OOXML), we first blow
# away the shape's text with an unrelated placeholder, then set it to
# the real final text. Paragraph-level formatting (e.g. algn="r") is
# preserved by the engine regardless.
#
# Note: the interpreter mis-parses a parenthesized concatenation
# expression when it's written directly in a function-call argument
# list, so every text value is built into a local variable first and
# only the variable is passed into the helper functions.

function Set-ShapeText($shape, [string]$text) {
    $tr = $shape.TextFrame.TextRange
    $tr.Text = "@@@PLACEHOLDER@@@"
    $tr2 = $shape.TextFrame.TextRange
    $tr2.Text = $text
}

function Get-ShapeById($slide, [int]$id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$NL = [char]13

# --- Food names (Nachos/Hotdog/Pretzels/Pickles/Pizza box) ---
# Before: Nachos, Hotdog, Pretzels, Pickles, Pizza
# After:  Pretzels , Nachos, Hotdog, Pizza, Pickles
$foodNames = Get-ShapeById $s 94
$foodNamesText = "Pretzels " + $NL + "Nachos" + $NL + "Hotdog" + $NL + "Pizza" + $NL + "Pickles"
Set-ShapeText $foodNames $foodNamesText

# --- Snacks names (Chips/Candy/Ring Pops/Suckers/Cotton Candy box) ---
# Before: Chips, Candy, Ring Pops, Suckers, Cotton Candy
# After:  Cotton Candy, Candy , Chips, Ring Pops, Suckers
$snackNames = Get-ShapeById $s 97
$snackNamesText = "Cotton Candy" + $NL + "Candy " + $NL + "Chips" + $NL + "Ring Pops" + $NL + "Suckers"
Set-ShapeText $snackNames $snackNamesText

# --- Drinks names (Pop/Water/Gatorade/Coffee/Bug Juice box) ---
# Before: Pop/Water, Gatorade, Coffee/Hot Chocolate, Bug Juice
# After:  Gatorade , Pop/Water, Coffee/Hot Chocolate, Bug Juice
$drinkNames = Get-ShapeById $s 98
$drinkNamesText = "Gatorade " + $NL + "Pop/Water" + $NL + "Coffee/Hot Chocolate" + $NL + "Bug Juice"
Set-ShapeText $drinkNames $drinkNamesText

# --- Food prices (matches food names box order) ---
# Before: $2.50, $2.50, $3.00, $1.00, $2.50
# After:  $3.00, $2.50, $2.50, $2.50, $1.00
$foodPrices = Get-ShapeById $s 7
$foodPricesText = "$3.00" + $NL + "$2.50" + $NL + "$2.50" + $NL + "$2.50" + $NL + "$1.00"
Set-ShapeText $foodPrices $foodPricesText

# --- Snacks prices (matches snack names box order) ---
# Before: $1.50, $2.00, $1.00, $.50, $4.00
# After:  $4.00, $2.00, $1.50, $1.00, $.50
$snackPrices = Get-ShapeById $s 9
$snackPricesText = "$4.00" + $NL + "$2.00" + $NL + "$1.50" + $NL + "$1.00" + $NL + "$.50"
Set-ShapeText $snackPrices $snackPricesText

# --- Drinks prices (matches drink names box order) ---
# Before: $2.00, $2.00/$3.00, $2.00, $2.00
# After:  $3.00, $2.00, $2.00, $2.00
$drinkPrices = Get-ShapeById $s 11
$drinkPricesText = "$3.00" + $NL + "$2.00" + $NL + "$2.00" + $NL + "$2.00"
Set-ShapeText $drinkPrices $drinkPricesText
